$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = 150
$ws.Range("C6").Value = 150
$ws.Range("C7").Value = 175
$ws.Range("C10").Value = 90
$ws.Range("C11").Value = 90
$ws.Range("C12").Value = 90
$ws.Range("C13").Value = 90
$ws.Range("C14").Value = 100
$ws.Range("C15").Value = 100
$ws.Range("C18").Value = 100

$ws.Activate()
$ws.Range("C1").Select()
